$wb = $excel.ActiveWorkbook

# --- Rename sheets: BF/CF (board feet / cubic feet) -> MBF/CCF (thousand board feet / hundred cubic feet) ---
$wsHarvest = $wb.Worksheets.Item("Harvest_BF")
$wsHarvest.Name = "Harvest_MBF"

$wsBfcf = $wb.Worksheets.Item("BFCF")
$wsBfcf.Name = "MBFCCF"

# --- HWP_MODEL_OPTIONS sheet updates ---
$wsOptions = $wb.Worksheets.Item("HWP_MODEL_OPTIONS")

# Unit conversion default changed from 1000 to 100 (M2)
$wsOptions.Range("M2").Value = 100

# Set explicit (best-fit) column widths for columns A:N so the header/value
# text is fully visible, matching the target "bestFit" widths.
$wsOptions.Columns.Item(1).ColumnWidth = 14.166666666666666
$wsOptions.Columns.Item(2).ColumnWidth = 7.833333333333333
$wsOptions.Columns.Item(3).ColumnWidth = 15.333333333333334
$wsOptions.Columns.Item(4).ColumnWidth = 14.666666666666666
$wsOptions.Columns.Item(5).ColumnWidth = 9.5
$wsOptions.Columns.Item(6).ColumnWidth = 14.666666666666666
$wsOptions.Columns.Item(7).ColumnWidth = 14.5
$wsOptions.Columns.Item(8).ColumnWidth = 29.833333333333332
$wsOptions.Columns.Item(9).ColumnWidth = 30
$wsOptions.Columns.Item(10).ColumnWidth = 7.5
$wsOptions.Columns.Item(11).ColumnWidth = 3.1666666666666665
$wsOptions.Columns.Item(12).ColumnWidth = 13.833333333333334
$wsOptions.Columns.Item(13).ColumnWidth = 5.833333333333333
$wsOptions.Columns.Item(14).ColumnWidth = 13.833333333333334

# Update selected cell on this sheet
$wsOptions.Activate() | Out-Null
$wsOptions.Range("Q29").Select() | Out-Null

# --- MBFCCF (formerly BFCF) sheet: update selected cell ---
$wsBfcf.Activate() | Out-Null
$wsBfcf.Range("I30").Select() | Out-Null

# Re-activate the sheet that was originally active/tabSelected
$wsOptions.Activate() | Out-Null
